$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("company_list")

# 2014/12
$ws.Range("D2").Value = 1102
$ws.Range("E2").Value = -118
$ws.Range("F2").Value = -118
$ws.Range("G2").Value = -80
$ws.Range("H2").Value = -75
$ws.Range("I2").Value = -73
$ws.Range("J2").Value = -2
$ws.Range("K2").Value = 997
$ws.Range("L2").Value = 655
$ws.Range("M2").Value = 342
$ws.Range("N2").Value = 341
$ws.Range("O2").Value = 2
$ws.Range("P2").Value = 203
$ws.Range("Q2").Value = 26
$ws.Range("R2").Value = 56
$ws.Range("S2").Value = -54
$ws.Range("T2").Value = 7
$ws.Range("U2").Value = 19
$ws.Range("V2").Value = 477
$ws.Range("W2").Value = -10.68
$ws.Range("X2").Value = -6.82
$ws.Range("Y2").Value = -19.54
$ws.Range("Z2").Value = -6.8
$ws.Range("AA2").Value = 191.25
$ws.Range("AB2").Value = 64.82
$ws.Range("AC2").Value = -1910
$ws.Range("AD2").Value = -3.94
$ws.Range("AE2").Value = 8395
$ws.Range("AF2").Value = 0.9
$ws.Range("AG2").Value = 0
$ws.Range("AH2").Value = 0
$ws.Range("AI2").Value = 0
$ws.Range("AJ2").Value = 4059609

# 2015/12
$ws.Range("D3").Value = 1185
$ws.Range("E3").Value = 23
$ws.Range("F3").Value = 23
$ws.Range("G3").Value = -1
$ws.Range("H3").Value = -7
$ws.Range("I3").Value = -4
$ws.Range("J3").Value = -3
$ws.Range("K3").Value = 936
$ws.Range("L3").Value = 600
$ws.Range("M3").Value = 336
$ws.Range("N3").Value = 337
$ws.Range("O3").Value = -1
$ws.Range("P3").Value = 203
$ws.Range("Q3").Value = -45
$ws.Range("R3").Value = -4
$ws.Range("S3").Value = -25
$ws.Range("T3").Value = 7
$ws.Range("U3").Value = -52
$ws.Range("V3").Value = 463
$ws.Range("W3").Value = 1.9
$ws.Range("X3").Value = -0.57
$ws.Range("Y3").Value = -1.11
$ws.Range("Z3").Value = -0.7
$ws.Range("AA3").Value = 178.69
$ws.Range("AB3").Value = 63.06
$ws.Range("AC3").Value = -93
$ws.Range("AD3").Value = -94.5
$ws.Range("AE3").Value = 8312
$ws.Range("AF3").Value = 1.06
$ws.Range("AG3").Value = 0
$ws.Range("AH3").Value = 0
$ws.Range("AI3").Value = 0
$ws.Range("AJ3").Value = 4059609

# 2016/12
$ws.Range("D4").Value = 1331
$ws.Range("E4").Value = 13
$ws.Range("F4").Value = 13
$ws.Range("G4").Value = 14
$ws.Range("H4").Value = 5
$ws.Range("I4").Value = 9
$ws.Range("J4").Value = -3
$ws.Range("K4").Value = 942
$ws.Range("L4").Value = 606
$ws.Range("M4").Value = 336
$ws.Range("N4").Value = 336
$ws.Range("O4").Value = 0
$ws.Range("P4").Value = 203
$ws.Range("Q4").Value = 26
$ws.Range("R4").Value = -29
$ws.Range("S4").Value = 7
$ws.Range("T4").Value = 51
$ws.Range("U4").Value = -25
$ws.Range("V4").Value = 463
$ws.Range("W4").Value = 1.01
$ws.Range("X4").Value = 0.41
$ws.Range("Y4").Value = 2.56
$ws.Range("Z4").Value = 0.58
$ws.Range("AA4").Value = 180.03
$ws.Range("AB4").Value = 66.81
$ws.Range("AC4").Value = 212
$ws.Range("AD4").Value = 45.65
$ws.Range("AE4").Value = 8284
$ws.Range("AF4").Value = 1.17
$ws.Range("AG4").Value = 0
$ws.Range("AH4").Value = 0
$ws.Range("AI4").Value = 0
$ws.Range("AJ4").Value = 4059609

# 2017/12
$ws.Range("D5").Value = 1617
$ws.Range("E5").Value = 66
$ws.Range("F5").Value = 66
$ws.Range("G5").Value = 65
$ws.Range("H5").Value = 51
$ws.Range("I5").Value = 54
$ws.Range("J5").Value = -2
$ws.Range("K5").Value = 1048
$ws.Range("L5").Value = 655
$ws.Range("M5").Value = 393
$ws.Range("N5").Value = 395
$ws.Range("O5").Value = -2
$ws.Range("P5").Value = 210
$ws.Range("Q5").Value = 145
$ws.Range("R5").Value = -11
$ws.Range("S5").Value = -10
$ws.Range("T5").Value = 8
$ws.Range("U5").Value = 137
$ws.Range("V5").Value = 433
$ws.Range("W5").Value = 4.09
$ws.Range("X5").Value = 3.16
$ws.Range("Y5").Value = 14.66
$ws.Range("Z5").Value = 5.14
$ws.Range("AA5").Value = 166.52
$ws.Range("AB5").Value = 91.61
$ws.Range("AC5").Value = 1318
$ws.Range("AD5").Value = 7.15
$ws.Range("AE5").Value = 9421
$ws.Range("AF5").Value = 1
$ws.Range("AG5").Value = 0
$ws.Range("AH5").Value = 0
$ws.Range("AI5").Value = 0
$ws.Range("AJ5").Value = 4197449

# 2018/12
$ws.Range("D6").Value = 1420
$ws.Range("E6").Value = 50
$ws.Range("F6").Value = 50
$ws.Range("G6").Value = 25
$ws.Range("H6").Value = 14
$ws.Range("I6").Value = 17
$ws.Range("K6").Value = 954
$ws.Range("L6").Value = 506
$ws.Range("M6").Value = 448
$ws.Range("N6").Value = 446
$ws.Range("P6").Value = 233
$ws.Range("Q6").Value = 4
$ws.Range("R6").Value = -7
$ws.Range("S6").Value = -43
$ws.Range("T6").Value = 4
$ws.Range("U6").Value = 0
$ws.Range("V6").Value = 343
$ws.Range("W6").Value = 3.52
$ws.Range("X6").Value = 0.97
$ws.Range("Y6").Value = 4.08
$ws.Range("Z6").Value = 1.37
$ws.Range("AA6").Value = 112.89
$ws.Range("AB6").Value = 95.73
$ws.Range("AC6").Value = 384
$ws.Range("AD6").Value = 23.45
$ws.Range("AE6").Value = 9581
$ws.Range("AF6").Value = 0.94
$ws.Range("AI6").Value = 0
$ws.Range("AJ6").Value = 4653805


# Row 6 (2018/12): AG6 and AH6 are removed entirely (no longer populated)
$ws.Range("AG6:AH6").ClearContents()

# Rows 7, 8, 9: all financial data (columns D through AJ) is removed, leaving
# only the A (index), B (company name) and C (period label) columns intact.
$ws.Range("D7:AJ7").ClearContents()
$ws.Range("D8:AJ8").ClearContents()
$ws.Range("D9:AJ9").ClearContents()
